# Insert a new weekly data row (row 324) in the Berenjena price sheet,
# shifting all existing rows 324-372 down by one (they become 325-373).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 324..372 down to 325..373 and make room for the new row 324.
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with the new weekly record.
$ws.Range("A324").Value = 6
$ws.Range("B324").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C324").Value = "Metropolitana"
$ws.Range("D324").Value = 45218
$ws.Range("E324").Value = 13
$ws.Range("F324").Value = 100112001
$ws.Range("G324").Value = "Berenjena"
$ws.Range("H324").Value = "Sin especificar"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 155
$ws.Range("K324").Value = 6000
$ws.Range("L324").Value = 7000
$ws.Range("M324").Value = 6645
$ws.Range("N324").Value = "$/caja 50 unidades"
$ws.Range("O324").Value = "Región de Arica y Parinacota"
$ws.Range("P324").Value = 133
$ws.Range("Q324").Value = 50
$ws.Range("R324").Value = "Hortaliza"
